$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.510.01'
$ws.Range("E2").Value = '  -2.27%  '

# Row 3
$ws.Range("D3").Value = '3.028.37'
$ws.Range("E3").Value = '  -1.41%  '

# Row 4
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.34'
$ws.Range("E5").Value = '  -0.55%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.54'
$ws.Range("E6").Value = '  +1.16%  '

# Row 7
$ws.Range("E7").Value = '  +0.12%  '

# Row 8
$ws.Range("D8").Value = '3.026.69'
$ws.Range("E8").Value = '  -1.28%  '

# Row 9
$ws.Range("E9").Value = '  +1.46%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.150'
$ws.Range("E10").Value = '  -2.68%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.15'
$ws.Range("E11").Value = '  -1.27%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").Value = '  -0.50%  '

# Row 13
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '34.73'
$ws.Range("E13").Value = '  +1.31%  '

# Row 14
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000222'
$ws.Range("E14").Value = '  -0.59%  '

# Row 15
$ws.Range("D15").Value = '3.523.69'

# Row 16
$ws.Range("E16").Value = '  +0.18%  '

# Row 17
$ws.Range("D17").Value = '61.662.10'
$ws.Range("E17").Value = '  -1.96%  '

# Row 18
$ws.Range("D18").Value = '3.036.08'
$ws.Range("E18").Value = '  -0.95%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.69'
$ws.Range("E19").Value = '  +0.71%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '471.80'
$ws.Range("E20").Value = '  -2.39%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.38'
$ws.Range("E21").Value = '  +0.39%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.687'
$ws.Range("E22").Value = '  -1.19%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.04'
$ws.Range("E23").Value = '  -1.36%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.00'
$ws.Range("E24").Value = '  +1.21%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.26'
$ws.Range("E25").Value = '  +1.12%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.70'
$ws.Range("E27").Value = '  +0.19%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.93'
$ws.Range("E28").Value = '  -2.47%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.07%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.93'
$ws.Range("E30").Value = '  +3.38%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.90'
$ws.Range("E31").Value = '  -0.48%  '

# Row 32
$ws.Range("E32").Value = '  +3.05%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.63'
$ws.Range("E33").Value = '  +4.67%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.31'
$ws.Range("E34").Value = '  -2.06%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.73'
$ws.Range("E35").Value = '  -1.83%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.97'
$ws.Range("E36").Value = '  -0.37%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '466.72'
$ws.Range("E37").Value = '  -2.80%  '

# Row 38
$ws.Range("D38").Value = '3.233.03'
$ws.Range("E38").Value = '  +4.67%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0798'
$ws.Range("E39").Value = '  +0.08%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0391'
$ws.Range("E40").Value = '  -0.70%  '

# Row 41
$ws.Range("E41").Value = '  +2.73%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.23'
$ws.Range("E42").Value = '  +1.45%  '

# Row 43
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '27.90'
$ws.Range("E43").Value = '  +13.45%  '

# Row 44
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.51'
$ws.Range("E44").Value = '  -5.19%  '

# Row 45
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.08%  '

# Row 46
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.249'
$ws.Range("E46").Value = '  -1.33%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.03'
$ws.Range("E47").Value = '  +0.59%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.109'
$ws.Range("E48").Value = '  +1.28%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '119.11'
$ws.Range("E49").Value = '  -1.78%  '

# Row 50
$ws.Range("D50").Value = '0.0₃0502'
$ws.Range("E50").Value = '  -7.35%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.04'
$ws.Range("E51").Value = '  +2.08%  '
